# Core Dataset refactor: the merge-diagnostic columns produced by
# date_proximity are now namespaced with an "_mp_" prefix, and the
# "_duplicates" indicator column has been dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O1").Value = "_mp_merge"
$ws.Range("P1").Value = "_mp_diff_days"
$ws.Range("Q1").Value = "_mp_abs_diff_days"

# "_duplicates" (column R) no longer exists - remove the whole column so the
# used range shrinks from A1:R13 back down to A1:Q13.
$ws.Columns.Item(18).Delete()

# The renamed headers are longer than the originals, so their columns grow
# to fit the new best-fit width.
$ws.Columns.Item(15).ColumnWidth = 9.5
$ws.Columns.Item(16).ColumnWidth = 11.17
$ws.Columns.Item(17).ColumnWidth = 14.5
